# Modifications pour utiliser XGBClassifier et ajuster les predictions
$wb = $excel.ActiveWorkbook

# --- Sheet "Valeurs reelles": rename the S+1/S+2/S+3 columns to the
#     classifier-output variants and replace the regression values in
#     columns C:E (rows 2-28) with the new class predictions. ---
$ws1 = $wb.Worksheets.Item("Valeurs réelles")

$ws1.Range("C1").Value = "PRIX EXP POMME GOLDEN FRANCE 136/200G CAT.I CAISSE_S+1_class"
$ws1.Range("D1").Value = "PRIX EXP POMME GOLDEN FRANCE 136/200G CAT.I CAISSE_S+2_class"
$ws1.Range("E1").Value = "PRIX EXP POMME GOLDEN FRANCE 136/200G CAT.I CAISSE_S+3_class"

$classValues = @(
    @(2, 0, 2, 2),
    @(3, 2, 2, 2),
    @(4, 2, 2, 2),
    @(5, 2, 2, 2),
    @(6, 2, 2, 2),
    @(7, 2, 2, 2),
    @(8, 2, 2, 2),
    @(9, 2, 2, 0),
    @(10, 2, 0, 1),
    @(11, 0, 1, 4),
    @(12, 1, 4, 0),
    @(13, 4, 0, 2),
    @(14, 0, 2, 4),
    @(15, 2, 4, 2),
    @(16, 4, 2, 0),
    @(17, 2, 0, 0),
    @(18, 0, 0, 4),
    @(19, 0, 4, 4),
    @(20, 4, 4, 3),
    @(21, 4, 3, 2),
    @(22, 3, 2, 3),
    @(23, 2, 3, 0),
    @(24, 3, 0, 2),
    @(25, 0, 2, 0),
    @(26, 2, 0, 2),
    @(27, 0, 2, 2),
    @(28, 2, 2, 2)
)

foreach ($entry in $classValues) {
    $r = $entry[0]
    $ws1.Cells.Item($r, 3).Value = $entry[1]
    $ws1.Cells.Item($r, 4).Value = $entry[2]
    $ws1.Cells.Item($r, 5).Value = $entry[3]
}

# --- Sheet "Predictions": the PRED_S1/S2/S3 regression outputs are all
#     zeroed out for rows 2-28 (the model now predicts classes elsewhere). ---
$ws2 = $wb.Worksheets.Item("Prédictions")

for ($r = 2; $r -le 28; $r++) {
    $ws2.Cells.Item($r, 2).Value = 0
    $ws2.Cells.Item($r, 3).Value = 0
    $ws2.Cells.Item($r, 4).Value = 0
}
